$d = $word.ActiveDocument

# The "Professional experience" table (5th table in the document) lists one
# "Tools" row per job/project entry. Each such row's value cell originally
# held template placeholder leftovers (commas and a stray "a" token). The
# client filled in some of these with real tool names; this change cleans
# each cell up to contain only the de-duplicated, comma-separated list of
# tools the client actually entered, and clears the cells the client left
# untouched (only template placeholder junk, no real tool) back to empty.

$t = $d.Tables.Item(5)

# Row 7  (",, a"              -> "abc")
$t.Cell(7,2).Range.Text = "abc"

# Row 15 (",, AWS, a"         -> "AWS")
$t.Cell(15,2).Range.Text = "AWS"

# Row 23 (",, a"              -> empty cell, no tool supplied by client)
[void]$t.Cell(23,2).Range.Paragraphs.Item(1).Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/></w:p>")

# Row 31 (",, Git"            -> "Git")
$t.Cell(31,2).Range.Text = "Git"

# Row 39 (",, Salesforce, a"  -> "Salesforce")
$t.Cell(39,2).Range.Text = "Salesforce"

# Row 47 (",, a"              -> empty cell, no tool supplied by client)
[void]$t.Cell(47,2).Range.Paragraphs.Item(1).Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/></w:p>")

# Row 55 (",, excel, a"       -> "excel")
$t.Cell(55,2).Range.Text = "excel"

# Row 63 ("C, C#, #, a, ,"    -> "C, C#")
$t.Cell(63,2).Range.Text = "C, C#"
